$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value2 = "nccFileAddress"
$ws.Range("B14").Value2 = $ws.Range("B2").Value2

$ws.Range("B14").Select()
